$p = $ppt.ActivePresentation

# --- Slide 7 (sldId 262, cId 1603588476): "Approximation:" content placeholder ---
$s7 = $p.Slides.Item(7)
$shp7 = $s7.Shapes.Item(4)   # "Content Placeholder 2" (id=11)

# Turn off autofit (was normAutofit w/ lnSpcReduction, now noAutofit)
$shp7.TextFrame.AutoSize = 0

# Add a new paragraph "{3, 4, 5}" after the existing "{0, 1, 2}" line
$shp7.TextFrame.TextRange.InsertAfter([char]13 + "{3, 4, 5}")

# --- Slide 8 (sldId 263, cId 1903942611): "Approximation:" content placeholder ---
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(4)   # "Content Placeholder 2" (id=11)

# Turn off autofit (was normAutofit w/ lnSpcReduction, now noAutofit)
$shp8.TextFrame.AutoSize = 0

# Add a new paragraph "{1, 3, 5, 6, 7, 8, 10}" after the existing "{1, 2, 5, 6, 7, 8, 9}" line
$tr8 = $shp8.TextFrame.TextRange
$beforeLen = $tr8.Length
$tr8.InsertAfter([char]13 + "{1, 3, 5, 6, 7, 8, 10}")

# The new paragraph text starts right after "beforeLen" characters + the inserted paragraph break
$newParaStart = $beforeLen + 2
$run1 = $tr8.Characters($newParaStart, 19)          # "{1, 3, 5, 6, 7, 8, "
$run2 = $tr8.Characters($newParaStart + 19, 3)       # "10}"

foreach ($r in @($run1, $run2)) {
    $r.Font.Bold = $false
    $r.Font.Italic = $false
    $r.Font.Underline = $false
    $r.Font.Strikethrough = $false
    $r.Font.Shadow = $false
    $r.Font.Name = "Calibri"
    $r.Font.Color.RGB = 0x3F3F3F
}
